# feat: add 2022-Q4 data
#
# The existing "2022-Q3" detail sheet is renamed to "2022-Q4" and gets the
# new quarter's fund-holding figures; a fresh sheet named "2022-Q3" is
# inserted right after it, carrying the figures the "2022-Q3" sheet used to
# hold. The "总计" (totals) summary sheet gains a new row for 2022-Q4 while
# keeping the old 2022-Q3 row.

$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item("总计")
$q3Sheet = $wb.Worksheets.Item("2022-Q3")

# ---------------------------------------------------------------------
# 1) Rename the original "2022-Q3" sheet to "2022-Q4" (keeps its
#    sheetId/r:id in place; the object reference $q3Sheet stays valid).
# ---------------------------------------------------------------------
$q3Sheet.Name = "2022-Q4"

# ---------------------------------------------------------------------
# 2) Insert a brand-new "2022-Q3" sheet right after it, pasting in the
#    formatting (borders/fonts/alignment) the detail sheets share.
#    (Column A on row 1 is never used on the detail sheets, so the header
#    row and data rows are copied separately to avoid touching A1.)
# ---------------------------------------------------------------------
$newQ3Sheet = $wb.Worksheets.Add($null, $q3Sheet)
$newQ3Sheet.Name = "2022-Q3"

$q3Sheet.Range("B1:H1").Copy()
$newQ3Sheet.Range("B1:H1").PasteSpecial(-4122)
$q3Sheet.Range("A2:H3").Copy()
$newQ3Sheet.Range("A2:H3").PasteSpecial(-4122)

$newQ3Sheet.Outline.SummaryBelow = $true
$newQ3Sheet.Outline.SummaryRight = $true

# ---------------------------------------------------------------------
# 3) Populate the brand-new "2022-Q3" sheet with the figures that used to
#    live on the "2022-Q3" sheet (moved data, unchanged values). A leading
#    apostrophe keeps numeric-looking text (leading zeros, "2.49", ...)
#    stored as text instead of being coerced to a number.
# ---------------------------------------------------------------------
$newQ3Sheet.Range("B1").Value = "基金代码"
$newQ3Sheet.Range("C1").Value = "基金名称"
$newQ3Sheet.Range("D1").Value = "基金规模"
$newQ3Sheet.Range("E1").Value = "股票总仓位"
$newQ3Sheet.Range("F1").Value = "仓位占比"
$newQ3Sheet.Range("G1").Value = "持有市值(亿元)"
$newQ3Sheet.Range("H1").Value = "仓位排名"

$newQ3Sheet.Range("A2").Value = 0
$newQ3Sheet.Range("B2").Value = "'014806"
$newQ3Sheet.Range("C2").Value = "国金量化精选混合C"
$newQ3Sheet.Range("D2").Value = "'2.49"
$newQ3Sheet.Range("E2").Value = "'80.89"
$newQ3Sheet.Range("F2").Value = "'0.74"
$newQ3Sheet.Range("G2").Value = "'0.0184"
$newQ3Sheet.Range("H2").Value = 7

$newQ3Sheet.Range("A3").Value = 1
$newQ3Sheet.Range("B3").Value = "'014805"
$newQ3Sheet.Range("C3").Value = "国金量化精选混合A"
$newQ3Sheet.Range("D3").Value = "'0.78"
$newQ3Sheet.Range("E3").Value = "'80.89"
$newQ3Sheet.Range("F3").Value = "'0.74"
$newQ3Sheet.Range("G3").Value = "'0.0058"
$newQ3Sheet.Range("H3").Value = 7

# ---------------------------------------------------------------------
# 4) Overwrite the renamed "2022-Q4" sheet with the new quarter's data
#    (only a single fund this quarter, so drop the now-unused 3rd row).
# ---------------------------------------------------------------------
$q4Sheet = $wb.Worksheets.Item("2022-Q4")

$q4Sheet.Range("B1").Value = "基金代码"
$q4Sheet.Range("C1").Value = "基金名称"
$q4Sheet.Range("D1").Value = "基金规模"
$q4Sheet.Range("E1").Value = "股票总仓位"
$q4Sheet.Range("F1").Value = "仓位占比"
$q4Sheet.Range("G1").Value = "持有市值(亿元)"
$q4Sheet.Range("H1").Value = "仓位排名"

$q4Sheet.Range("A2").Value = 0
$q4Sheet.Range("B2").Value = "'001397"
$q4Sheet.Range("C2").Value = "建信精工制造指数增强"
$q4Sheet.Range("D2").Value = "'0.47"
$q4Sheet.Range("E2").Value = "'91.21"
$q4Sheet.Range("F2").Value = "'1.93"
$q4Sheet.Range("G2").Value = "'0.0091"
$q4Sheet.Range("H2").Value = 8

$q4Sheet.Rows(3).Delete()

# ---------------------------------------------------------------------
# 5) Update the "总计" (totals) sheet: row 2 becomes 2022-Q4, and a new
#    row 3 (copied formatting from row 2) keeps the old 2022-Q3 totals.
# ---------------------------------------------------------------------
$summary.Range("A2:D2").Copy()
$summary.Range("A3:D3").PasteSpecial(-4122)

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 2
$summary.Range("D3").Value = 0.02

$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0.01

# Keep "总计" as the active sheet/selection, matching the original workbook.
$summary.Activate()
$summary.Range("A1").Select() | Out-Null
